# repull data, push all data, mean calculation
# Update the dSF (column F) values on specific rows to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    5  = -2
    10 = -1
    11 = -5
    14 = -1
    17 = 1
    22 = 0
    23 = 1
    29 = -1
    35 = -5
    36 = 9
    40 = 0
    41 = 0
    42 = -1
    45 = -1
    46 = -2
    49 = -3
    51 = -2
    52 = 0
    53 = -4
    61 = -3
    62 = -4
    68 = -5
    73 = -2
    74 = -6
    75 = -2
    78 = -6
    79 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
